$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "G2"
$ws.Range("B3").Value = "Test1jq"
$ws.Range("C3").Value = "Daily"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 45860
$ws.Range("F3").Value = 30

# Match the date style used in E2 (numFmtId 165 -> YYYY-MM-DD)
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
